# wms_kaart_database.xlsx - column-name casing fix + related view/formatting touch-ups
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Blad1
$ws2 = $wb.Worksheets.Item(2)   # Blad2
$ws3 = $wb.Worksheets.Item(3)   # Blad3

# ---------------------------------------------------------------------------
# 1) Rename the "VIEW ATTRIBUTEN" values so the real column names are spelled
#    out (snake_case instead of the old run-together abbreviations).
# ---------------------------------------------------------------------------
$ws1.Range("F2:F4").Value = "identificatie_lokaalid, bgt_functie, plus_functie"
$ws1.Range("F5:F7").Value = "identificatie_lokaalid, bgt_type, plus_type"

# ---------------------------------------------------------------------------
# 2) Add the next generation of the (auto)filter defined name, continuing the
#    existing _FilterDatabase_0 / _FilterDatabase_0_0 sequence.
# ---------------------------------------------------------------------------
$ws1.Names.Add("_xlnm._FilterDatabase_0_0_0", "=Blad1!`$A`$1:`$J`$4")

# ---------------------------------------------------------------------------
# 3) Widen the columns (workbook was re-saved through Excel, which re-measured
#    the "best fit" column widths against its own font metrics).
# ---------------------------------------------------------------------------
$offset = 5/6

$sheet1Widths = @(
    17.0242914979757,
    39.9919028340081,
    21.4817813765182,
    20.4534412955466,
    41.246963562753,
    44.4331983805668,
    40.5627530364373,
    23.1943319838057,
    18.165991902834,
    18.6234817813765,
    9.1417004048583
)
for ($i = 0; $i -lt $sheet1Widths.Length; $i++) {
    $ws1.Columns.Item($i + 1).ColumnWidth = $sheet1Widths[$i] - $offset
}

$ws2.Columns.Item(1).ColumnWidth = 8.79757085020243 - $offset
$ws3.Columns.Item(1).ColumnWidth = 8.79757085020243 - $offset

# ---------------------------------------------------------------------------
# 4) Move the Blad1 view / selection: the unfrozen (top) pane now scrolls to
#    column B and both panes' selections move accordingly.
# ---------------------------------------------------------------------------
$aw1 = $excel.Windows.Item(1)
$aw1.ScrollColumn = 2
$aw1.ScrollRow = 1

$topSel = $ws1.Range("B1")
$bottomArea1 = $ws1.Range("F5:F7")
$bottomArea2 = $ws1.Range("E11")
$excel.Union($bottomArea1, $bottomArea2).Select()
$topSel.Select()

# Blad2 / Blad3 selections also change to include F5:F7 plus A1 (A1 active).
$excel.ActiveWindow.ActivateNext() | Out-Null

foreach ($ws in @($ws2, $ws3)) {
    $ws.Activate()
    $a1 = $ws.Range("A1")
    $f5f7 = $ws.Range("F5:F7")
    $excel.Union($f5f7, $a1).Select()
}

$ws1.Activate()
